# TableSxx_GO_KEGG.xlsx — "Show SYMBOl instead of ENTREZID"
#
# The KEGG enrichment sheet's geneID column (column I) held slash-separated
# NCBI ENTREZID lists (e.g. "10000/5894/2033"). This replaces them with the
# equivalent gene SYMBOL lists (e.g. "AKT3/RAF1/EP300"), and nudges the
# saved cell-selection on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "GO enrichment"
$ws2 = $wb.Worksheets.Item(2)   # "KEGG enrichment"

# ENTREZID -> SYMBOL for every gene referenced in the KEGG geneID column.
$entrez2symbol = @{
    "10000"  = "AKT3"
    "5894"   = "RAF1"
    "2033"   = "EP300"
    "57492"  = "ARID1B"
    "6310"   = "ATXN1"
    "2776"   = "GNAQ"
    "23389"  = "MED13L"
}

function Convert-GeneIds($value) {
    $ids = $value -split "/"
    $symbols = @()
    foreach ($id in $ids) {
        $symbols += $entrez2symbol[$id]
    }
    return [string]::Join("/", $symbols)
}

# Column I ("geneID") runs from row 2 through row 61 on the KEGG sheet.
$lastRow = $ws2.Cells.Item($ws2.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 61) { $lastRow = 61 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws2.Cells.Item($r, 9)
    $cur = $cell.Value2
    if ($cur -ne $null -and $cur -ne "") {
        $cell.Value = Convert-GeneIds $cur
    }
}

# Saved cell-selection moves (cosmetic, matches the authored commit).
$ws2.Range("L25").Select()
$ws1.Activate()
$ws1.Range("Z32").Select()
